$d = $word.ActiveDocument
$maxIter = 20
$i = 0
while ($d.Revisions.Count -gt 0 -and $i -lt $maxIter) {
    $r = $d.Revisions.Item(1)
    Write-Output "count=$($d.Revisions.Count) type=$($r.Type) text=$($r.Range.Text)"
    $r.Accept()
    $i = $i + 1
}
